$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 106143380.984658
$ws.Range("D2").Value = 48.509731

$ws.Range("B3").Value = 15530776.874209
$ws.Range("D3").Value = 3.548944
$ws.Range("E3").Value = 0.03039

$ws.Range("B4").Value = 487942795.675565
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = -632.019392
$ws.Range("H5").Value = -1278.047091
$ws.Range("I5").Value = 14.008308
$ws.Range("J5").Value = 0.056696

$ws.Range("G6").Value = -191.00834
$ws.Range("H6").Value = -886.643385
$ws.Range("I6").Value = 504.626705
$ws.Range("J6").Value = 0.793767

$ws.Range("G7").Value = 441.011052
$ws.Range("H7").Value = -84.896277
$ws.Range("I7").Value = 966.91838
$ws.Range("J7").Value = 0.119883
